$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell 2 4 '43.085.27'
Set-TextCell 2 5 '  +3.01%  '

Set-TextCell 3 4 '2.299.41'
Set-TextCell 3 5 '  +1.97%  '

Set-TextCell 4 5 '  +0.10%  '

Set-TextCell 5 4 '310.32'
Set-TextCell 5 5 '  +1.98%  '

Set-TextCell 6 4 '100.90'
Set-TextCell 6 5 '  +6.29%  '

Set-TextCell 7 5 '  +2.47%  '

Set-TextCell 8 5 '  +0.07%  '

Set-TextCell 9 4 '0.523'
Set-TextCell 9 5 '  +7.48%  '

Set-TextCell 10 4 '36.15'
Set-TextCell 10 5 '  +4.37%  '

Set-TextCell 11 4 '0.0824'
Set-TextCell 11 5 '  +4.80%  '

Set-TextCell 12 5 '  +0.78%  '

Set-TextCell 13 5 '  +7.64%  '

Set-TextCell 14 4 '2.657.51'
Set-TextCell 14 5 '  +2.17%  '

Set-TextCell 15 4 '15.00'
Set-TextCell 15 5 '  +4.96%  '

Set-TextCell 16 4 '2.306.32'
Set-TextCell 16 5 '  +2.35%  '

Set-TextCell 17 4 '0.807'
Set-TextCell 17 5 '  +2.52%  '

Set-TextCell 18 4 '43.034.41'
Set-TextCell 18 5 '  +3.17%  '

Set-TextCell 19 4 '12.51'
Set-TextCell 19 5 '  +1.67%  '

Set-TextCell 20 4 '0.0₃0922'
Set-TextCell 20 5 '  +2.73%  '

Set-TextCell 21 4 '6.07'
Set-TextCell 21 5 '  +2.26%  '

Set-TextCell 22 4 '68.47'
Set-TextCell 22 5 '  +0.94%  '

Set-TextCell 23 4 '240.07'
Set-TextCell 23 5 '  +1.54%  '

Set-TextCell 24 4 '2.01'
Set-TextCell 24 5 '  +4.94%  '

Set-TextCell 25 4 '2.61'
Set-TextCell 25 5 '  +1.84%  '

Set-TextCell 26 4 '0.999'
Set-TextCell 26 5 '  -0.02%  '

Set-TextCell 27 4 '24.63'
Set-TextCell 27 5 '  +4.53%  '

Set-TextCell 28 4 '38.46'
Set-TextCell 28 5 '  +6.17%  '

Set-TextCell 29 4 '9.65'
Set-TextCell 29 5 '  +2.29%  '

Set-TextCell 30 4 '2.11'
Set-TextCell 30 5 '  +0.25%  '

Set-TextCell 31 4 '167.87'
Set-TextCell 31 5 '  +5.25%  '

Set-TextCell 32 4 '5.32'
Set-TextCell 32 5 '  +2.49%  '

Set-TextCell 33 5 '  +0.13%  '

Set-TextCell 34 5 '  -0.72%  '

Set-TextCell 35 4 '17.72'
Set-TextCell 35 5 '  +4.04%  '

Set-TextCell 36 4 '0.0739'
Set-TextCell 36 5 '  +0.91%  '

Set-TextCell 37 5 '  +2.82%  '

Set-TextCell 38 5 '  +0.45%  '

$ws.Cells.Item(39, 2).Value = 'Stellar'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 39 4 '0.116'
Set-TextCell 39 5 '  +2.06%  '

$ws.Cells.Item(40, 2).Value = 'ARBITRUM'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 40 4 '1.83'
Set-TextCell 40 5 '  +1.12%  '

Set-TextCell 41 5 '  +5.78%  '

Set-TextCell 42 4 '2.30'
Set-TextCell 42 5 '  -2.58%  '

$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 43 4 '1.972.65'
Set-TextCell 43 5 '  +0.07%  '

$ws.Cells.Item(44, 2).Value = 'VeChain'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 44 4 '0.0288'
Set-TextCell 44 5 '  +2.39%  '

$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 45 4 '19.15'
Set-TextCell 45 5 '  +2.19%  '

Set-TextCell 46 5 '  +4.17%  '

Set-TextCell 47 4 '9.83'
Set-TextCell 47 5 '  +0.06%  '

Set-TextCell 48 4 '55.60'
Set-TextCell 48 5 '  +4.90%  '

Set-TextCell 49 5 '  +15.68%  '

Set-TextCell 50 4 '2.526.50'
Set-TextCell 50 5 '  +1.97%  '

Set-TextCell 51 5 '  +2.57%  '
